$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / shared-string renames ---
$ws.Cells.Item(1,3).Value  = "GDP"                    # Gross_National_Income -> GDP
$ws.Cells.Item(1,5).Value  = "Budget_Previous_Year"   # NGO_Country_Budget_Previous_Year -> Budget_Previous_Year
$ws.Cells.Item(1,6).Value  = "LatinAmerica"            # Vision_ONGD_LatinAmerica -> LatinAmerica
$ws.Cells.Item(1,7).Value  = "Africa"                  # Vision_ONGD_Africa -> Africa
$ws.Cells.Item(1,8).Value  = "Confessional"            # Vision_ONGD_Confessional -> Confessional
$ws.Cells.Item(1,9).Value  = "Universal"               # Vision_ONGD_Universal -> Universal
$ws.Cells.Item(1,32).Value = "Donor_Aid_Budget"        # Total_subvencion_en_el_Pais_y_Anyo -> Donor_Aid_Budget
$ws.Cells.Item(1,33).Value = "Total_Funds"             # Total_Fondos -> Total_Funds
$ws.Cells.Item(1,34).Value = "%_Private_Funds"         # Proporcion_Fondos_Privados -> %_Private_Funds
$ws.Cells.Item(1,35).Value = "%_MAE_Funds"             # Proporcion_Fondos_MAE -> %_MAE_Funds
$ws.Cells.Item(1,39).Value = "Delegation"              # Delegacion -> Delegation

# --- Column C (GDP) value updates, rows 2-42 ---
$ws.Cells.Item(2,3).Value  = 2771.04675450926
$ws.Cells.Item(3,3).Value  = 2870.311589353206
$ws.Cells.Item(4,3).Value  = 1460.056109840828
$ws.Cells.Item(5,3).Value  = 9502.243585046588
$ws.Cells.Item(6,3).Value  = 4547.50930098406
$ws.Cells.Item(7,3).Value  = 2100.656463590606
$ws.Cells.Item(8,3).Value  = 19868.07076233724
$ws.Cells.Item(9,3).Value  = 17288.8595992193
$ws.Cells.Item(10,3).Value = 951.6879611168786
$ws.Cells.Item(11,3).Value = 665.6274194933962
$ws.Cells.Item(12,3).Value = 1503.870423231357
$ws.Cells.Item(13,3).Value = 10385.96443195552
$ws.Cells.Item(14,3).Value = 4633.590358399045
$ws.Cells.Item(15,3).Value = 1357.563719132622
$ws.Cells.Item(16,3).Value = 492.3430015592067
$ws.Cells.Item(17,3).Value = 17610.30663334184
$ws.Cells.Item(18,3).Value = 982.980837581714
$ws.Cells.Item(19,3).Value = 2965.153206179127
$ws.Cells.Item(20,3).Value = 691.8942672110555
$ws.Cells.Item(21,3).Value = 1577.487171555845
$ws.Cells.Item(22,3).Value = 4921.848409120176
$ws.Cells.Item(23,3).Value = 5360.226632400601
$ws.Cells.Item(24,3).Value = 18254.09644617555
$ws.Cells.Item(25,3).Value = 1000.829216794104
$ws.Cells.Item(26,3).Value = 5122.180090208862
$ws.Cells.Item(27,3).Value = 16764.42871195103
$ws.Cells.Item(28,3).Value = 1032.277326842402
$ws.Cells.Item(29,3).Value = 5295.682695961288
$ws.Cells.Item(30,3).Value = 3252.634165082374
$ws.Cells.Item(31,3).Value = 1640.18070024053
$ws.Cells.Item(32,3).Value = 1060.095015975378
$ws.Cells.Item(33,3).Value = 711.3043470146426
$ws.Cells.Item(34,3).Value = 846.386841468855
$ws.Cells.Item(35,3).Value = 3314.741082534716
$ws.Cells.Item(36,3).Value = 1751.664428859304
$ws.Cells.Item(37,3).Value = 1093.134170274031
$ws.Cells.Item(38,3).Value = 731.9993357350996
$ws.Cells.Item(39,3).Value = 871.998368594318
$ws.Cells.Item(40,3).Value = 12358.30403621203
$ws.Cells.Item(41,3).Value = 729.8559996981501
$ws.Cells.Item(42,3).Value = 729.6614300490079

# Row 43, column C: was text ".." -> becomes numeric 0
$ws.Cells.Item(43,3).Value = 0

# --- Column AL (Colony) flips from 0 to 1 on rows 11, 20, 41 ---
$ws.Cells.Item(11,38).Value = 1
$ws.Cells.Item(20,38).Value = 1
$ws.Cells.Item(41,38).Value = 1
